$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "wHzCP787"
$ws.Range("B2").Value = 23102867
$ws.Range("C2").Value = "qmlznnr19"
$ws.Range("D2").Value = "Wb9%G!m3"
$ws.Range("F2").Value = "jOcWgytr"
$ws.Range("G2").Value = "Xcmq"

$ws.Range("A3").Value = "GbKze769"
$ws.Range("B3").Value = 23102866
$ws.Range("C3").Value = "hwuazkk64"
$ws.Range("D3").Value = "q29#VQ!u"
$ws.Range("F3").Value = "zbJvSuxA"
$ws.Range("G3").Value = "EvAS"

$ws.Range("A4").Value = "wyJnJ700"
$ws.Range("B4").Value = 23102865
$ws.Range("C4").Value = "cbsegfy45"
$ws.Range("D4").Value = "c92Z%!bK"
$ws.Range("F4").Value = "gfPlNKBQ"
$ws.Range("G4").Value = "tmsq"

$ws.Range("A5").Value = "DKgWk194"
$ws.Range("B5").Value = 23102864
$ws.Range("C5").Value = "rxcurlr46"
$ws.Range("D5").Value = "T%D4km9!"
$ws.Range("F5").Value = "ArBpbFMc"
$ws.Range("G5").Value = "oAgk"

$ws.Range("A6").Value = "kNImS133"
$ws.Range("B6").Value = 23102863
$ws.Range("C6").Value = "qraweuw95"
$ws.Range("D6").Value = "Qs42&#Nx"
$ws.Range("F6").Value = "VjRRxrXN"
$ws.Range("G6").Value = "ueoo"

$ws.Range("A7").Value = "StcrL389"
$ws.Range("B7").Value = 23102862
$ws.Range("C7").Value = "nvomkpz46"
$ws.Range("D7").Value = "N5u!7P%n"
$ws.Range("F7").Value = "jMIPSJoG"
$ws.Range("G7").Value = "PcRk"

$ws.Range("A8").Value = "XwAza117"
$ws.Range("B8").Value = 23102861
$ws.Range("C8").Value = "chrrhgm45"
$ws.Range("D8").Value = "V&7a#yT2"
$ws.Range("F8").Value = "XfnVKMys"
$ws.Range("G8").Value = "aYhX"

$ws.Range("A9").Value = "sVQmP200"
$ws.Range("B9").Value = 23102860
$ws.Range("C9").Value = "gyrjiem33"
$ws.Range("D9").Value = "m`$&y2GH9"
$ws.Range("F9").Value = "XpEgOxTT"
$ws.Range("G9").Value = "mjGZ"

$ws.Range("A10").Value = "ewbVq834"
$ws.Range("B10").Value = 23102859
$ws.Range("C10").Value = "wdhfgky90"
$ws.Range("D10").Value = "Gf&m9#3S"
$ws.Range("F10").Value = "JTvfXRsY"
$ws.Range("G10").Value = "qyAm"

$ws.Range("A11").Value = "JjqAC825"
$ws.Range("B11").Value = 23102858
$ws.Range("C11").Value = "rfcwakg29"
$ws.Range("D11").Value = "U28r`$W&a"
$ws.Range("F11").Value = "kKtemfbI"
$ws.Range("G11").Value = "qHsu"

$ws.Range("A12").Value = "MCYbS352"
$ws.Range("B12").Value = 23102857
$ws.Range("C12").Value = "zlztkft69"
$ws.Range("D12").Value = "Wt8#5N&m"
$ws.Range("F12").Value = "uRuUfYJL"
$ws.Range("G12").Value = "auie"

$ws.Range("A13").Value = "PkaGP273"
$ws.Range("B13").Value = 23102856
$ws.Range("C13").Value = "nnljccf74"
$ws.Range("D13").Value = "bT`$6%Dw9"
$ws.Range("F13").Value = "iUuCoqrC"
$ws.Range("G13").Value = "CtMa"

$ws.Range("A14").Value = "HWAnk167"
$ws.Range("B14").Value = 23102855
$ws.Range("C14").Value = "cbxhkth25"
$ws.Range("D14").Value = "Pn2%4f&B"
$ws.Range("F14").Value = "jsoVqMpX"
$ws.Range("G14").Value = "SSCw"

# Remove row 15 (former iAuthor TC row no longer present)
$ws.Rows("15:15").Delete()

